# Generate Report for Archive
#
# 1) Every "Ready for handoff" status cell becomes "In Translation".
# 2) The "Status"-ish columns (zh-cn/de-de on the Overview sheet, Status on
#    the per-locale sheets) get narrower now that the status text is shorter.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $v = $cell.Value2
            if (($v -is [string]) -and ($v -eq "Ready for handoff")) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# Overview sheet: columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

# zh-cn / de-de sheets: column C (Status)
$wb.Worksheets.Item("zh-cn").Range("C1").ColumnWidth = 12.5
$wb.Worksheets.Item("de-de").Range("C1").ColumnWidth = 12.5
